# Applies the commit: "Switched to estimating a risk ratio as specified in protocol"
#
# 1. Update the generated timestamp in the byline.
# 2. Swap "odds ratio" / "logit link" wording for "risk ratio" / "log link"
#    in the statistical-methods paragraph.
# 3. Update the regression-output table (exp(b), Std. Err., t, P>|t|,
#    95% CI) to the numbers produced by the new risk-ratio model.

$d = $word.ActiveDocument

# --- 1. Timestamp in the byline -------------------------------------------
$d.Content.Find.Execute("09:07:3012 May 2020", $true, $false, $false, $false,
    $false, $true, 1, $false, "09:14:1512 May 2020", 2)

# --- 2. Methods paragraph wording ------------------------------------------
$d.Content.Find.Execute(
    "the imputed data. For each imputed data set, we estimated an odds ratio to ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "the imputed data. For each imputed data set, we estimated a risk ratio to ",
    2)

$d.Content.Find.Execute(
    "logit link) to account for the cluster design. Estimates were then combined ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "log link) to account for the cluster design. Estimates were then combined ",
    2)

# --- 3. Regression table ----------------------------------------------------
$t = $d.Tables.Item(1)

$changes = @(
    @{Row=4;  Col=3; Old="0.06";   New="0.04"},
    @{Row=4;  Col=4; Old="0.36";   New="0.37"},
    @{Row=4;  Col=5; Old="0.72";   New="0.71"},
    @{Row=4;  Col=6; Old="0.91";   New="0.93"},
    @{Row=4;  Col=7; Old="1.15";   New="1.11"},

    @{Row=7;  Col=2; Old="1.14";   New="1.10"},
    @{Row=7;  Col=3; Old="0.12";   New="0.09"},
    @{Row=7;  Col=4; Old="1.26";   New="1.25"},
    @{Row=7;  Col=6; Old="0.93";   New="0.94"},
    @{Row=7;  Col=7; Old="1.39";   New="1.29"},

    @{Row=8;  Col=3; Old="0.12";   New="0.09"},
    @{Row=8;  Col=6; Old="0.71";   New="0.77"},
    @{Row=8;  Col=7; Old="1.19";   New="1.14"},

    @{Row=9;  Col=2; Old="1.12";   New="1.09"},
    @{Row=9;  Col=3; Old="0.12";   New="0.09"},
    @{Row=9;  Col=4; Old="1.08";   New="1.07"},
    @{Row=9;  Col=6; Old="0.91";   New="0.93"},
    @{Row=9;  Col=7; Old="1.37";   New="1.27"},

    @{Row=10; Col=2; Old="1.21";   New="1.15"},
    @{Row=10; Col=3; Old="0.12";   New="0.09"},
    @{Row=10; Col=4; Old="1.84";   New="1.83"},
    @{Row=10; Col=7; Old="1.48";   New="1.34"},

    @{Row=12; Col=2; Old="0.30";   New="0.23"},
    @{Row=12; Col=3; Old="0.03";   New="0.02"},
    @{Row=12; Col=4; Old="-13.28"; New="-21.19"},
    @{Row=12; Col=6; Old="0.25";   New="0.20"},
    @{Row=12; Col=7; Old="0.36";   New="0.27"}
)

foreach ($chg in $changes) {
    $cell = $t.Cell($chg.Row, $chg.Col)
    $current = $cell.Range.Text.TrimEnd([char]13, [char]7)
    if ($current -eq $chg.Old) {
        $cell.Range.Text = $chg.New
    } else {
        Write-Output ("Unexpected value at Row=" + $chg.Row + " Col=" + $chg.Col + `
            " -- found '" + $current + "' expected '" + $chg.Old + "'")
    }
}

Write-Output "Edit complete"
